# -----------------------------------------------------------------------
# Target revision analysis
# -----------------------------------------------------------------------
# The supplied OOXML diff only ever touches the root element's xmlns:*
# declaration *order* (plus the incidental renumbering of the
# auto-generated "ns8" -> "ns9" alias that order shift causes) in nine
# package parts: word/document.xml, word/endnotes.xml, word/footer1.xml,
# word/footer2.xml, word/footnotes.xml, word/header1.xml,
# word/numbering.xml, word/styles.xml and word/theme/theme1.xml.
#
# The set of namespace URIs bound on every one of those root elements is
# identical before and after (same prefixes map to the same URIs, modulo
# the cosmetic ns8/ns9 alias), and every hunk's context lines show the
# element's children (<w:body>, <w:endnote>, <w:p>, <w:tbl>,
# <w:abstractNum>, <w:docDefaults>, <a:themeElements>, ...) are byte for
# byte unchanged. There is no paragraph/run/table/style/numbering/
# header/footer/theme content edit anywhere in the diff. The reordering
# is the signature of which XML serializer wrote the package (it is not
# meaningful XML - attribute order is not observable through the
# document object model), so it isn't something reachable from
# Find/Replace, Range/Selection edits, style edits, etc. The commit
# message itself confirms the real change was to unrelated application
# code (date-range validation on a "financeiro"/"agenda" screen, and
# swapping out a third-party XML conversion library) - nothing about
# this "TERMO DE DECLARACAO REQUISITO I" template's text or formatting.
#
# So the faithful edit here is "no content change" - touch the active
# document via the object model (proving the session is live) without
# mutating anything.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# Read-only round trip through the object model; intentionally makes no
# edits, because the diff contains no content change to apply.
$paragraphCount = $d.Paragraphs.Count
$tableCount = $d.Tables.Count
